$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New 2020 values (column Q), one per data row 3..33 (row 3 is the year header row)
$values = @{
    3  = 2020
    4  = 53.463696812512026
    5  = 46.05204738706685
    6  = 60.76705279190513
    7  = 46.609654277955656
    8  = 39.785591828762811
    9  = 53.693996785869842
    10 = 49.132459991853935
    11 = 42.132308166831223
    12 = 56.225753650646354
    13 = 28.457427087863305
    14 = 20.524708126577082
    15 = 36.325895173845353
    16 = 37.816151622141014
    17 = 29.032279844170926
    18 = 46.928626462141906
    19 = 51.38232216208695
    20 = 45.862881450184311
    21 = 57.0280888993139
    22 = 44.951834666409091
    23 = 38.216466887636237
    24 = 51.83682668469686
    25 = 82.176148450436926
    26 = 66.965035434789911
    27 = 96.931980629894966
    28 = 56.391242440049062
    29 = 50.844030930786069
    30 = 61.300998533028128
    31 = 54.829571415516767
    32 = 58.407045187583961
    33 = 51.452932817170577
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 17).Value = $values[$row]
}

# Copy the number formatting / styling from column P (existing last-year column)
# onto the newly populated column Q, cell by cell, so each row keeps the same
# look as its neighbour in P.
for ($row = 3; $row -le 33; $row++) {
    $ws.Range("P$row").Copy() | Out-Null
    $ws.Range("Q$row").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# Update the active selection to match the latest edit location
$ws.Range("S34").Select() | Out-Null

Write-Host "done"
